$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "choices" sheet: add the two new "choice_filter" helper columns
#    (filter1 / filter2) used by the new select_one/select_multiple
#    choice-filter questions, plus their numeric values per choice row.
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

$choices.Range("D1").Value = "filter1"
$choices.Range("E1").Value = "filter2"

$choices.Range("D6").Value = 1
$choices.Range("E6").Value = 1

$choices.Range("D7").Value = 1
$choices.Range("E7").Value = 2

$choices.Range("D8").Value = 1
$choices.Range("E8").Value = 3

$choices.Range("D9").Value = 2
$choices.Range("E9").Value = 4

$choices.Range("D10").Value = 2
$choices.Range("E10").Value = 5

# ---------------------------------------------------------------------------
# 2. "survey" sheet: tidy up a couple of stray / duplicate cell styles left
#    over in column C (rows 3-13 previously used near-duplicate, unused
#    style variants) so they match the tidy styles used elsewhere on the
#    sheet.
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

$survey.Range("B3").Copy() | Out-Null
$survey.Range("C3:C4").PasteSpecial(-4122) | Out-Null

$survey.Range("A5").Copy() | Out-Null
$survey.Range("C5:C13").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3. Switch the active/selected sheet from "survey" to "choices", landing the
#    selection on B1 (matches the workbook being saved while viewing the
#    choices sheet).
# ---------------------------------------------------------------------------
$choices.Select() | Out-Null
$choices.Range("B1").Select() | Out-Null
